# Completed Connecting Excel data to Webops Class
# Replaces the placeholder "Details"/"SoftwareDeveloper"/"DataScientist"
# sample data with the real keyword/contact data used by the scraping class.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Details": personal/contact info instead of the old Title/URL row
# ---------------------------------------------------------------------
$wsDetails = $wb.Worksheets.Item("Details")

$wsDetails.Range("A2").Value = "first_name"
$wsDetails.Range("B2").Value = "Rahul"

$wsDetails.Range("A3").Value = "last_name"
$wsDetails.Range("B3").Value = "Doe"

$wsDetails.Range("A4").Value = "phone"
$wsDetails.Range("B4").Value = "999-999-9999"

$wsDetails.Range("A5").Value = "email"
$wsDetails.Range("B5").Value = "test@gmail.com"

$wsDetails.Range("A6").Value = "city"
$wsDetails.Range("B6").Value = "testCity"

# the old B2 cell ("https://www.google.com/") carried a fill style that is
# no longer used anywhere in the workbook - drop it so the cell goes back
# to the default "Normal" style
$wsDetails.Range("B2").ClearFormats()

# ---------------------------------------------------------------------
# Sheet "SoftwareDeveloper": search-keyword breakdown
# ---------------------------------------------------------------------
$wsSoftware = $wb.Worksheets.Item("SoftwareDeveloper")

$wsSoftware.Range("A2").Value = "exactPhrase"
$wsSoftware.Range("B2").Value = "software"

$wsSoftware.Range("A3").Value = "Keywords"
$wsSoftware.Range("B3").Value = "developer engineer"

$wsSoftware.Range("A4").Value = "Experience"
$wsSoftware.Range("B4").Value = 3

$wsSoftware.Range("A5").Value = "skills"
$wsSoftware.Range("B5").Value = "Java; Python"

# ---------------------------------------------------------------------
# Sheet "DataScientist": search-keyword breakdown
# ---------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("DataScientist")

$wsData.Range("A2").Value = "exactPhrase"
$wsData.Range("B2").Value = "data"

$wsData.Range("A3").Value = "Keywords"
$wsData.Range("B3").Value = "analysist scientist"

$wsData.Range("A4").Value = "Experience"
$wsData.Range("B4").Value = 3

$wsData.Range("A5").Value = "skills"
$wsData.Range("B5").Value = "SQL; Python"

# ---------------------------------------------------------------------
# Selections / active sheet.  Select on the two non-active sheets first,
# then select+activate "Details" last so it ends up the tab that is
# actually shown when the workbook is reopened.
# ---------------------------------------------------------------------
$wsSoftware.Range("A2:B5").Select() | Out-Null
$wsData.Range("A6").Select() | Out-Null
$wsDetails.Range("A10").Select() | Out-Null
